# Apply the tracked changes:
#  1. Remove the _GoBack bookmark from its old location (end of the
#     paragraph containing "ذی‌نفعان").
#  2. Strip the stray <w:rtl/> from the paragraph-mark run properties of
#     the paragraph containing "معیارهای توسعه".
#  3. Insert three new numbered list paragraphs right after that
#     paragraph, and re-create the _GoBack bookmark at the end of the
#     third one.

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark -----------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- locate the target paragraph ("معیارهای توسعه") ------------------------
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match [regex]::Escape("معیارهای توسعه")) {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)

# --- 2. Rewrite this paragraph without the paragraph-mark <w:rtl/> --------
$pXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="17009052" w14:textId="77777777" w:rsidR="00A04D1E" w:rsidRPr="002316D0" w:rsidRDefault="00A04D1E" w:rsidP="00A04D1E" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:bidi/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r w:rsidRPr="002316D0"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>معیارهای توسعه</w:t></w:r></w:p>'
$p.Range.InsertXML($pXml)

# --- 3. Insert three new paragraphs after it -------------------------------
$p.Range.InsertParagraphAfter()
$pA = $d.Paragraphs.Item($targetIndex + 1)
$pA.Range.InsertParagraphAfter()
$pB = $d.Paragraphs.Item($targetIndex + 2)
$pB.Range.InsertParagraphAfter()
$pC = $d.Paragraphs.Item($targetIndex + 3)

$xmlA = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:bidi/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">1-کمینه هزینه جهت استقرار </w:t></w:r></w:p>'
$pA.Range.InsertXML($xmlA)

$xmlB = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:bidi/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">2- </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">تسریع در استقرار سیستم </w:t></w:r></w:p>'
$pB.Range.InsertXML($xmlB)

# The _GoBack bookmark is re-created here, embedded directly in the OOXML
# payload (Bookmarks.Add with a collapsed range sitting exactly on a
# paragraph mark is unreliable in this runtime), right after the "3-" run
# and before the paragraph mark - matching its original position relative
# to the last run in its paragraph.
$xmlC = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:bidi/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="44"/><w:szCs w:val="44"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>3-</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$pC.Range.InsertXML($xmlC)

Write-Output "done"
